$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 23

# Columns A and D look like numbers/dates to Excel's auto-detection, so force
# them to be stored as text (matching the other rows in this sheet) by
# temporarily applying a text number format, then resetting the style back
# to Normal so no new cell style is introduced.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-01-13"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "22:15:33"
$ws.Cells.Item($row, 3).Value = "Monday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "02"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 126875
$ws.Cells.Item($row, 6).Value = 143483
$ws.Cells.Item($row, 7).Value = 169306
$ws.Cells.Item($row, 8).Value = 152060
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142736
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193047
$ws.Cells.Item($row, 14).Value = 115436
$ws.Cells.Item($row, 15).Value = 45881
$ws.Cells.Item($row, 16).Value = 28503
$ws.Cells.Item($row, 17).Value = 65496
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 48941
$ws.Cells.Item($row, 20).Value = -1
